$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update Road (row 3) stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 264
$wsOff.Range("C3").Value = 161
$wsOff.Range("F3").Value = 8

# Sheet "DEF" - update Road (row 3) stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 281
$wsDef.Range("C3").Value = 205
$wsDef.Range("D3").Value = 52
$wsDef.Range("E3").Value = 21
$wsDef.Range("G3").Value = 3
